# Update "F" column (想去人数 / interest count) values to match latest scrape
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(6, 6).Value = 443  # was 442
$ws.Cells.Item(8, 6).Value = 2008  # was 2001
$ws.Cells.Item(10, 6).Value = 36  # was 35
$ws.Cells.Item(11, 6).Value = 33  # was 32
$ws.Cells.Item(14, 6).Value = 1324  # was 1323
$ws.Cells.Item(19, 6).Value = 8  # was 7
$ws.Cells.Item(20, 6).Value = 451  # was 448
$ws.Cells.Item(22, 6).Value = 145  # was 142
$ws.Cells.Item(23, 6).Value = 7036  # was 7032
$ws.Cells.Item(24, 6).Value = 7036  # was 7032
$ws.Cells.Item(25, 6).Value = 7618  # was 7605
$ws.Cells.Item(28, 6).Value = 179  # was 178
$ws.Cells.Item(33, 6).Value = 125  # was 50
$ws.Cells.Item(36, 6).Value = 39  # was 38
$ws.Cells.Item(38, 6).Value = 1385  # was 1380
$ws.Cells.Item(39, 6).Value = 15  # was 14
$ws.Cells.Item(41, 6).Value = 280  # was 278
$ws.Cells.Item(42, 6).Value = 685  # was 684
$ws.Cells.Item(46, 6).Value = 218  # was 213
$ws.Cells.Item(48, 6).Value = 81  # was 80
$ws.Cells.Item(49, 6).Value = 134  # was 129

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(3, 6).Value = 25  # was 24
$ws.Cells.Item(5, 6).Value = 51  # was 50
$ws.Cells.Item(17, 6).Value = 282  # was 281

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(3, 6).Value = 2569  # was 2565
$ws.Cells.Item(4, 6).Value = 261  # was 259
$ws.Cells.Item(5, 6).Value = 123  # was 121

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(6, 6).Value = 261  # was 259
$ws.Cells.Item(7, 6).Value = 123  # was 121
$ws.Cells.Item(10, 6).Value = 25  # was 24
$ws.Cells.Item(11, 6).Value = 443  # was 442
$ws.Cells.Item(12, 6).Value = 2008  # was 2001
$ws.Cells.Item(13, 6).Value = 36  # was 35
$ws.Cells.Item(14, 6).Value = 33  # was 32
$ws.Cells.Item(17, 6).Value = 1324  # was 1323
$ws.Cells.Item(19, 6).Value = 8  # was 7
$ws.Cells.Item(20, 6).Value = 451  # was 448
$ws.Cells.Item(22, 6).Value = 145  # was 142
$ws.Cells.Item(23, 6).Value = 51  # was 50
$ws.Cells.Item(24, 6).Value = 7036  # was 7032
$ws.Cells.Item(25, 6).Value = 7036  # was 7032
$ws.Cells.Item(26, 6).Value = 7618  # was 7605
$ws.Cells.Item(33, 6).Value = 39  # was 38
$ws.Cells.Item(34, 6).Value = 1385  # was 1380
$ws.Cells.Item(35, 6).Value = 15  # was 14
$ws.Cells.Item(38, 6).Value = 280  # was 278
$ws.Cells.Item(41, 6).Value = 685  # was 684
$ws.Cells.Item(46, 6).Value = 218  # was 213
$ws.Cells.Item(47, 6).Value = 81  # was 80
$ws.Cells.Item(48, 6).Value = 134  # was 129
$ws.Cells.Item(50, 6).Value = 282  # was 281
